# Update sheet names and corresponding KPI title cells (A1) to use
# zero-padded two-digit numbering (e.g. GP1 -> GP01, BP1 -> BP01, ...).

$wb = $excel.ActiveWorkbook

$renames = @(
    @{ Index = 1; OldName = "GP1"; NewName = "GP01"; OldTitle = "KPI GP1 - Global Perf 1"; NewTitle = "KPI GP01 - Global Perf 1" },
    @{ Index = 2; OldName = "GP2"; NewName = "GP02"; OldTitle = "KPI GP2 - Global Perf 2"; NewTitle = "KPI GP02 - Global Perf 2" },
    @{ Index = 3; OldName = "BP1"; NewName = "BP01"; OldTitle = "KPI BP1 - Business Process 1"; NewTitle = "KPI BP01 - Business Process 1" },
    @{ Index = 4; OldName = "BP2"; NewName = "BP02"; OldTitle = "KPI BP2 - Business Process 2"; NewTitle = "KPI BP02 - Business Process 2" },
    @{ Index = 5; OldName = "BP3"; NewName = "BP03"; OldTitle = "KPI BP3 - Business Process 3"; NewTitle = "KPI BP03 - Business Process 3" },
    @{ Index = 6; OldName = "BP4"; NewName = "BP04"; OldTitle = "KPI BP4 - Business Process 4"; NewTitle = "KPI BP04 - Business Process 4" },
    @{ Index = 7; OldName = "BP5"; NewName = "BP05"; OldTitle = "KPI BP5 - Business Process 5"; NewTitle = "KPI BP05 - Business Process 5" },
    @{ Index = 8; OldName = "BP6"; NewName = "BP06"; OldTitle = "KPI BP6 - Business Process 6"; NewTitle = "KPI BP06 - Business Process 6" },
    @{ Index = 9; OldName = "BP7"; NewName = "BP07"; OldTitle = "KPI BP7 - Business Process 7"; NewTitle = "KPI BP07 - Business Process 7" }
)

foreach ($item in $renames) {
    $ws = $wb.Worksheets.Item($item.Index)
    $ws.Range("A1").Value = $item.NewTitle
    $ws.Name = $item.NewName
}
